$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new diary entry row (row 23).
# Set values in the same order the original author's shared-string table
# was grown: date (A), content (C), time (B), then the META note (F).
$ws.Range("A23").Value = "4 marras"
$ws.Range("C23").Value = "Kangassimulaation numeerisen epästabiiliuden selvittely, "
$ws.Range("B23").Value = "18.00-18.45"
$ws.Range("F23").Value = "Merkkivirhe söi miestä liian monta tuntia tässä vaiheessa tutkintoa."

# Carry over the same cell formatting used by neighboring rows instead of
# minting new styles: B -> time format (like B11/B18), C/F -> wrapped text
# (like C22/F22).
$ws.Range("B11").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("F22").Copy()
$ws.Range("F23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New row wraps onto two lines, same as row 21.
$ws.Rows.Item(23).RowHeight = 29

# Fix the window anchoring: scroll the view down one row and move the
# active-cell/selection flag off the old last row (F22) onto the new one (G23).
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("G23").Select() | Out-Null
